$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.022.98'
$ws.Range('D3').Value = '1.672.60'
$ws.Range('E3').Value = '  +2.89%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('E6').Value = '  +1.98%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0617'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.10'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.86%  '
$ws.Range('E11').Value = '  +4.62%  '
$ws.Range('D12').Value = '1.909.23'
$ws.Range('E12').Value = '  +2.73%  '
$ws.Range('D13').Value = '1.671.15'
$ws.Range('E13').Value = '  +3.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.08'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '65.81'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.83%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.520'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '27.050.05'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '234.75'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.71'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.45'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.34%  '
$ws.Range('E23').Value = '  +1.87%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.25'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.28'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.90'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.33'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').Value = '1.447.99'
$ws.Range('E33').Value = '  -5.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.13'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.09%  '
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.569'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.892'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.72%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0169'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +11.62%  '
$ws.Range('E43').Value = '  +2.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '65.70'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.15%  '
$ws.Range('D45').Value = '1.818.63'
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.780'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.33'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +1.59%  '
$ws.Range('E49').Value = '  +4.16%  '
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.66'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.68%  '
